$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# --- Version/Change/Author/Date table updates ---
Replace-Text "1.0" "1.2.5"
Replace-Text "Creation" "Update"
Replace-Text "Fabrício Araújo" "Julio Paiva"
Replace-Text "09/07/2020" "31/05/2023"

# --- Precondition table text fix ---
Replace-Text "O usuario devidamente autenticado e na tela inicial do sistema" "O usuário devidamente autenticado e na tela inicial do sistema."

# --- Main flow step 1 wording update ---
Replace-Text "1. Beneficiário O usuario acessa o caso de uso atraves do menu. af[1,2,3,4]" "1. Beneficiário Acessa o caso de uso através do menu. af[1,2,3,4]"

# --- Alternative flow step 2 wording updates (trailing period) ---
Replace-Text "2. System Apresenta a tela de Detalhar Diárias " "2. System Apresenta a tela de Detalhar Diárias. "
Replace-Text "2. System Apresenta a tela de Analisar Prestação de Contas " "2. System Apresenta a tela de Analisar Prestação de Contas. "
Replace-Text "2. System Apresenta a tela de Cancelar Solicitação de Diária " "2. System Apresenta a tela de Cancelar Solicitação de Diária. "
